$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (match original inlineStr formatting)
$textCells = @("D5", "D6", "D7", "D10", "D11", "D12", "D13", "D15", "D17", "D18", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D42", "D44", "D48", "D50")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated values from the cryptos data refresh
$ws.Range("D2").Value = "69.816.42"
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("D3").Value = "3.516.66"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "616.46"
$ws.Range("E5").Value = "  +5.47%  "
$ws.Range("D6").Value = "192.35"
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  +0.73%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -3.68%  "
$ws.Range("D10").Value = "0.665"
$ws.Range("E10").Value = "  +2.00%  "
$ws.Range("D11").Value = "53.47"
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("D12").Value = "0.0000309"
$ws.Range("E12").Value = "  -3.56%  "
$ws.Range("D13").Value = "9.60"
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").Value = "4.084.52"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "620.06"
$ws.Range("E15").Value = "  +9.40%  "
$ws.Range("D16").Value = "69.849.34"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("D17").Value = "19.02"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "12.66"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("D19").Value = "3.509.84"
$ws.Range("E19").Value = "  -2.04%  "
$ws.Range("D22").Value = "108.79"
$ws.Range("E22").Value = "  +15.44%  "
$ws.Range("D23").Value = "17.20"
$ws.Range("E23").Value = "  -3.87%  "
$ws.Range("D24").Value = "4.71"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("D26").Value = "3.11"
$ws.Range("E26").Value = "  +5.98%  "
$ws.Range("D27").Value = "11.00"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("D28").Value = "9.72"
$ws.Range("E28").Value = "  +4.33%  "
$ws.Range("D29").Value = "34.36"
$ws.Range("E29").Value = "  +5.57%  "
$ws.Range("D30").Value = "6.98"
$ws.Range("E30").Value = "  -3.38%  "
$ws.Range("E31").Value = "  +1.54%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.116"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D33").Value = "3.91"
$ws.Range("E33").Value = "  +3.76%  "
$ws.Range("D34").Value = "63.43"
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("D35").Value = "3.11"
$ws.Range("E35").Value = "  -4.96%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.663.05"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "518.43"
$ws.Range("E38").Value = "  -2.76%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "3.65"
$ws.Range("E39").Value = "  +5.82%  "
$ws.Range("D40").Value = "0.393"
$ws.Range("E40").Value = "  -4.70%  "
$ws.Range("D41").Value = "0.0₃0780"
$ws.Range("E41").Value = "  -2.60%  "
$ws.Range("D42").Value = "36.72"
$ws.Range("E42").Value = "  -4.51%  "
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").Value = "0.0471"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("E46").Value = "  +3.14%  "
$ws.Range("E47").Value = "  -3.76%  "
$ws.Range("D48").Value = "8.78"
$ws.Range("E48").Value = "  -5.56%  "
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("D50").Value = "132.18"
$ws.Range("E50").Value = "  -2.12%  "
$ws.Range("E51").Value = "  -4.84%  "
